$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "-"
$ws.Range("C3").Value = "['MEC-3B-M. Motor Endot.', -, -, -]"
$ws.Range("E3").Value = "-"
$ws.Range("C4").Value = "['MEC-3B-M. Motor Endot.', -, -, -]"
$ws.Range("E4").Value = "-"
$ws.Range("C6").Value = "['MEC-3B-M. Motor Endot.', -, -, -]"
$ws.Range("B8").Value = "[-, 'MEC-3B-M. Motor Endot.', -, -]"
$ws.Range("E8").Value = "-"
$ws.Range("F10").Value = "-"
$ws.Range("F11").Value = "[-, -, -, 'MEC-3A-M. Motor Endot.']"
$ws.Range("F12").Value = "[-, -, -, 'MEC-3A-M. Motor Endot.']"
$ws.Range("F14").Value = "[-, -, -, 'MEC-3A-M. Motor Endot.']"
$ws.Range("E16").Value = "[-, -, 'MEC-3A-M. Motor Endot.', -]"
$ws.Range("F16").Value = "-"
$ws.Range("C18").Value = "[-, -, -, 'MEC-1NA-Manut. Mot. End.']"
$ws.Range("D18").Value = "['MEC-1NB-Manut. Mot. End.', -, -, -]"
$ws.Range("E18").Value = "-"
$ws.Range("C19").Value = "[-, -, -, 'MEC-1NA-Manut. Mot. End.']"
$ws.Range("E19").Value = "-"
$ws.Range("C20").Value = "[-, -, -, 'MEC-1NA-Manut. Mot. End.']"
$ws.Range("E20").Value = "-"
$ws.Range("C21").Value = "[-, -, -, 'MEC-1NA-Manut. Mot. End.']"
$ws.Range("E21").Value = "-"
